$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '48.197.25'
$ws.Range("E2").Value = '  +0.91%  '
$ws.Range("D3").Value = '2.502.93'
$ws.Range("E3").Value = '  +0.21%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '320.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.41%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '107.66'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.45%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.527'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.66%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.539'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.82%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.74'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.66%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.25'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +8.62%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0816'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.59%  '
$ws.Range("E13").Value = '  +0.00%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.15'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.92%  '
$ws.Range("D15").Value = '2.897.30'
$ws.Range("E15").Value = '  +0.43%  '
$ws.Range("D16").Value = '2.511.50'
$ws.Range("E16").Value = '  +0.60%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.840'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.92%  '
$ws.Range("D18").Value = '48.062.05'
$ws.Range("E18").Value = '  +0.96%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.09'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.69%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.76'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.72%  '
$ws.Range("D21").Value = '0.0₃0941'
$ws.Range("E21").Value = '  -0.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.74'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '278.52'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +12.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.02'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.87%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.54'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.54%  '
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.71'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.42%  '
$ws.Range("E28").Value = '  -0.19%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.76'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.52%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.140'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.18'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.06%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.48'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.85%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.72'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.62%  '
$ws.Range("B34").Value = 'FirstDigitalUSD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.01'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.09%  '
$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.34'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0780'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.02%  '
$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.95'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.78%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.66'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.56%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.91'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.111'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.57%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '121.94'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.80%  '
$ws.Range("E42").Value = '  +0.02%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.36'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.00%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0302'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.33%  '
$ws.Range("D45").Value = '2.015.06'
$ws.Range("E45").Value = '  +0.98%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.18'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.57%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.00'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.31%  '
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.85'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.29%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.02'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.27%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.17'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '80.38'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.37%  '
